# Auto-generated edit script: updates "想去人数" (F column) counts
# across sheets "展览", "演出", and "全部类型" per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 81
$ws.Range("F3").Value = 122
$ws.Range("F5").Value = 347
$ws.Range("F8").Value = 11518
$ws.Range("F12").Value = 2082
$ws.Range("F13").Value = 896
$ws.Range("F14").Value = 36
$ws.Range("F16").Value = 226
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 1175
$ws.Range("F19").Value = 150
$ws.Range("F20").Value = 236
$ws.Range("F21").Value = 723
$ws.Range("F22").Value = 132
$ws.Range("F23").Value = 250
$ws.Range("F24").Value = 2394
$ws.Range("F25").Value = 710
$ws.Range("F26").Value = 3406
$ws.Range("F27").Value = 1041
$ws.Range("F28").Value = 773
$ws.Range("F32").Value = 962
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 50
$ws.Range("F35").Value = 249
$ws.Range("F36").Value = 9
$ws.Range("F38").Value = 1718
$ws.Range("F39").Value = 4342
$ws.Range("F40").Value = 5433
$ws.Range("F42").Value = 109
$ws.Range("F43").Value = 25
$ws.Range("F44").Value = 145
$ws.Range("F45").Value = 248
$ws.Range("F47").Value = 21
$ws.Range("F48").Value = 4088
$ws.Range("F49").Value = 91

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4134
$ws.Range("F9").Value = 105
$ws.Range("F11").Value = 545
$ws.Range("F15").Value = 1
$ws.Range("F19").Value = 76

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 81
$ws.Range("F7").Value = 347
$ws.Range("F9").Value = 11518
$ws.Range("F12").Value = 2082
$ws.Range("F13").Value = 36
$ws.Range("F15").Value = 226
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 1175
$ws.Range("F18").Value = 150
$ws.Range("F19").Value = 236
$ws.Range("F20").Value = 4134
$ws.Range("F22").Value = 723
$ws.Range("F23").Value = 132
$ws.Range("F24").Value = 250
$ws.Range("F25").Value = 710
$ws.Range("F26").Value = 1041
$ws.Range("F29").Value = 773
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 9
$ws.Range("F37").Value = 109
$ws.Range("F38").Value = 145
$ws.Range("F39").Value = 248
$ws.Range("F42").Value = 21
$ws.Range("F43").Value = 4088
$ws.Range("F45").Value = 76
$ws.Range("F48").Value = 91
